# Refresh generated dataset values in column V (Pontuacao) for each round.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("V2").Value = 63.76
$ws.Range("V3").Value = 49.36
$ws.Range("V4").Value = 66.86
$ws.Range("V5").Value = 65.7
$ws.Range("V6").Value = 54.66
$ws.Range("V8").Value = 63.76
$ws.Range("V9").Value = 43.56
$ws.Range("V10").Value = 61.96
$ws.Range("V11").Value = 25.16
$ws.Range("V12").Value = 57.26
$ws.Range("V13").Value = 65.06
$ws.Range("V14").Value = 38.66
$ws.Range("V15").Value = 71.45999999999999
$ws.Range("V16").Value = 46.79
$ws.Range("V17").Value = 81.76000000000001
$ws.Range("V18").Value = 54.95
$ws.Range("V19").Value = 50.76
$ws.Range("V20").Value = 56.09
$ws.Range("V22").Value = 71.36
$ws.Range("V23").Value = 59.69
$ws.Range("V24").Value = 48.29
$ws.Range("V25").Value = 50.69
$ws.Range("V27").Value = 48.5
$ws.Range("V28").Value = 59.36
$ws.Range("V30").Value = 58.26
$ws.Range("V31").Value = 71.16
$ws.Range("V32").Value = 61.96
$ws.Range("V33").Value = 43.56
